$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2432.6667
$ws.Range("I29").Value = 2199
$ws.Range("J29").Value = 2549.5
$ws.Range("K29").Value = 6597
$ws.Range("L29").Value = 7648.5
$ws.Range("M29").Value = -6316
$ws.Range("N29").Value = -8210.5

$ws.Range("H33").Value = 155.85715
$ws.Range("I33").Value = 155.85715
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 155.85715
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 73.14285000000001
$ws.Range("N33").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H58").Value = 1845.8
$ws.Range("I58").Value = 882.25
$ws.Range("J58").Value = 5700
$ws.Range("K58").Value = 2646.75
$ws.Range("L58").Value = 17100
$ws.Range("M58").Value = -2496.75
$ws.Range("N58").Value = -17400

$ws.Range("H64").Value = 6000
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -5752
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 6000
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -5142
$ws.Range("N67").ClearContents()

$ws.Range("H74").Value = 2965.5
$ws.Range("I74").Value = 2965.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2965.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2029.5

$ws.Range("H77").Value = 2965.5
$ws.Range("I77").Value = 2965.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 14827.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -10147.5

$ws.Range("H87").Value = 80000.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 80000.5
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 80000.5
$ws.Range("N87").Value = -82496.5

$ws.Range("H90").Value = 80000.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 80000.5
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 240001.5
$ws.Range("N90").Value = -252481.5

$ws.Range("H103").Value = 213.33333
$ws.Range("I103").Value = 213.33333
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 639.99999
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -53.99999000000003
$ws.Range("N103").ClearContents()

$ws.Range("H113").Value = 2498.3333
$ws.Range("I113").Value = 2747.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2747.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 506.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11744.24
$ws.Range("I32").Value = 7558.8945
$ws.Range("J32").Value = 24997.834
$ws.Range("K32").Value = 7558.8945
$ws.Range("L32").Value = 24997.834
$ws.Range("M32").Value = -7271.8945
$ws.Range("N32").Value = -25571.834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 61397
$ws.Range("I62").Value = 4963.3335
$ws.Range("J62").Value = 399999
$ws.Range("K62").Value = 4963.3335
$ws.Range("L62").Value = 399999
$ws.Range("M62").Value = -4339.3335
$ws.Range("N62").Value = -401247

$ws.Range("H65").Value = 61397
$ws.Range("I65").Value = 4963.3335
$ws.Range("J65").Value = 399999
$ws.Range("K65").Value = 24816.6675
$ws.Range("L65").Value = 1999995
$ws.Range("M65").Value = -21696.6675
$ws.Range("N65").Value = -2006235

$ws.Range("H99").Value = 15853.869
$ws.Range("I99").Value = 14280
$ws.Range("J99").Value = 16693.268
$ws.Range("K99").Value = 14280
$ws.Range("L99").Value = 16693.268
$ws.Range("M99").Value = -12782
$ws.Range("N99").Value = -19689.268

$ws.Range("H122").Value = 2050.4443
$ws.Range("I122").Value = 2092.48
$ws.Range("J122").Value = 1525
$ws.Range("K122").Value = 6277.440000000001
$ws.Range("L122").Value = 4575
$ws.Range("M122").Value = -3827.440000000001

$ws.Range("H126").Value = 15853.869
$ws.Range("I126").Value = 14280
$ws.Range("J126").Value = 16693.268
$ws.Range("K126").Value = 42840
$ws.Range("L126").Value = 50079.804
$ws.Range("M126").Value = -40370
$ws.Range("N126").Value = -55019.804

$ws.Range("H134").Value = 3169.2942
$ws.Range("I134").Value = 2408.1667
$ws.Range("J134").Value = 4996
$ws.Range("K134").Value = 7224.500100000001
$ws.Range("L134").Value = 14988
$ws.Range("M134").Value = -4689.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -731
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -798
$ws.Range("N27").ClearContents()

$ws.Range("H114").Value = 1486.875
$ws.Range("I114").Value = 985
$ws.Range("J114").Value = 1654.1666
$ws.Range("K114").Value = 2955
$ws.Range("L114").Value = 4962.4998
$ws.Range("M114").Value = 299
$ws.Range("N114").Value = -11470.4998

$ws.Range("H132").Value = 1919.8
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 2199.6667
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 19797.0003
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -24857.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3922.2666
$ws.Range("I132").Value = 2400
$ws.Range("J132").Value = 6966.8
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 20900.4
$ws.Range("M132").Value = -4670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2289.6924
$ws.Range("I7").Value = 2363.4443
$ws.Range("J7").Value = 2123.75
$ws.Range("K7").Value = 2363.4443
$ws.Range("L7").Value = 2123.75
$ws.Range("M7").Value = -2251.4443

$ws.Range("H43").Value = 5000000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5000000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 5000000
$ws.Range("N43").Value = -5000386

$ws.Range("H126").Value = 2289.6924
$ws.Range("I126").Value = 2363.4443
$ws.Range("J126").Value = 2123.75
$ws.Range("K126").Value = 7090.3329
$ws.Range("L126").Value = 6371.25
$ws.Range("M126").Value = -4620.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 19655
$ws.Range("I30").Value = 19500
$ws.Range("J30").Value = 19810
$ws.Range("K30").Value = 19500
$ws.Range("L30").Value = 19810
$ws.Range("M30").Value = -19393
$ws.Range("N30").Value = -20024

$ws.Range("H64").Value = 57249.75
$ws.Range("I64").Value = 47999.5
$ws.Range("J64").Value = 66500
$ws.Range("K64").Value = 47999.5
$ws.Range("L64").Value = 66500
$ws.Range("M64").Value = -47751.5

$ws.Range("H67").Value = 57249.75
$ws.Range("I67").Value = 47999.5
$ws.Range("J67").Value = 66500
$ws.Range("K67").Value = 47999.5
$ws.Range("L67").Value = 66500
$ws.Range("M67").Value = -47141.5

$ws.Range("H126").Value = 2001.8334
$ws.Range("I126").Value = 1335.2222
$ws.Range("J126").Value = 4001.6667
$ws.Range("K126").Value = 4005.6666
$ws.Range("L126").Value = 12005.0001
$ws.Range("M126").Value = -1535.6666
